# Update crypto price/volume data per Dec 22 2023 GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.849.47'
$ws.Range("E2").Value = '  -0.68%  '

$ws.Range("D3").Value = '2.310.55'
$ws.Range("E3").Value = '  +2.54%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '94.87'
$c.ClearFormats()
$ws.Range("E5").Value = '  +8.11%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '268.62'
$c.ClearFormats()
$ws.Range("E6").Value = '  -1.41%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.620'
$c.ClearFormats()
$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +1.79%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '44.94'
$c.ClearFormats()
$ws.Range("E10").Value = '  -1.30%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0942'
$c.ClearFormats()
$ws.Range("E11").Value = '  +1.35%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '8.16'
$c.ClearFormats()
$ws.Range("E12").Value = '  +6.14%  '

$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").Value = '2.650.26'
$ws.Range("E14").Value = '  +2.19%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.43'
$c.ClearFormats()
$ws.Range("E15").Value = '  +2.61%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.867'
$c.ClearFormats()
$ws.Range("E16").Value = '  +8.98%  '

$ws.Range("D17").Value = '2.319.05'
$ws.Range("E17").Value = '  +2.17%  '

$ws.Range("D18").Value = '43.755.45'
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("E19").Value = '  +3.93%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.37'
$c.ClearFormats()
$ws.Range("E20").Value = '  +5.91%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '71.41'
$c.ClearFormats()
$ws.Range("E21").Value = '  +1.24%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '238.07'
$c.ClearFormats()
$ws.Range("E22").Value = '  +1.56%  '

$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.ClearFormats()
$ws.Range("E23").Value = '  -4.48%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.66'
$c.ClearFormats()
$ws.Range("E24").Value = '  +8.64%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  +4.78%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.51'
$c.ClearFormats()
$ws.Range("E27").Value = '  -1.31%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.48'
$c.ClearFormats()
$ws.Range("E28").Value = '  -1.92%  '

$ws.Range("E29").Value = '  -1.10%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '38.49'
$c.ClearFormats()
$ws.Range("E30").Value = '  -3.78%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '22.42'
$c.ClearFormats()
$ws.Range("E31").Value = '  +7.15%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '172.01'
$c.ClearFormats()
$ws.Range("E32").Value = '  -1.80%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0898'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.10%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.52'
$c.ClearFormats()
$ws.Range("E34").Value = '  +2.22%  '

$ws.Range("E35").Value = '  +1.74%  '

$ws.Range("E36").Value = '  -1.75%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0358'
$c.ClearFormats()
$ws.Range("E37").Value = '  +1.63%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.42'
$c.ClearFormats()
$ws.Range("E38").Value = '  +0.83%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.42'
$c.ClearFormats()
$ws.Range("E39").Value = '  -2.12%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.235'
$c.ClearFormats()
$ws.Range("E40").Value = '  +14.73%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.32'
$c.ClearFormats()
$ws.Range("E41").Value = '  +4.95%  '

$ws.Range("E42").Value = '  +19.53%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '12.19'
$c.ClearFormats()
$ws.Range("E43").Value = '  -4.13%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.46'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.00%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '62.18'
$c.ClearFormats()
$ws.Range("E45").Value = '  -4.20%  '

$ws.Range("E46").Value = '  +6.01%  '

$ws.Range("E47").Value = '  +3.43%  '

$ws.Range("E48").Value = '  -0.33%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.22'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.92%  '

$ws.Range("D50").Value = '2.531.62'
$ws.Range("E50").Value = '  +2.19%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.422'
$c.ClearFormats()
$ws.Range("E51").Value = '  -2.06%  '

Write-Output "Applied cryptos update"